$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 34412.5
$ws.Range("I28").Value = 41095.2
$ws.Range("K28").Value = 41095.2
$ws.Range("M28").Value = -40610.2
$ws.Range("H43").Value = 6688.6665
$ws.Range("I43").Value = 3475.25
$ws.Range("J43").Value = 8295.375
$ws.Range("K43").Value = 3475.25
$ws.Range("L43").Value = 8295.375
$ws.Range("M43").Value = -3406.25
$ws.Range("N43").Value = -8433.375
$ws.Range("H63").Value = 75000
$ws.Range("J63").Value = 75000
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76248
$ws.Range("H66").Value = 75000
$ws.Range("J66").Value = 75000
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -231240
$ws.Range("H86").Value = 3243.3
$ws.Range("I86").Value = 3299.2
$ws.Range("J86").Value = 3187.4
$ws.Range("K86").Value = 3299.2
$ws.Range("L86").Value = 3187.4
$ws.Range("M86").Value = -2176.2
$ws.Range("N86").Value = -5433.4
$ws.Range("H89").Value = 3243.3
$ws.Range("I89").Value = 3299.2
$ws.Range("J89").Value = 3187.4
$ws.Range("K89").Value = 16496
$ws.Range("L89").Value = 15937
$ws.Range("M89").Value = -10880
$ws.Range("N89").Value = -27169
$ws.Range("H105").Value = 72499.5
$ws.Range("J105").Value = 72499.5
$ws.Range("L105").Value = 72499.5
$ws.Range("N105").Value = -79487.5
$ws.Range("H132").Value = 1234.6222
$ws.Range("J132").Value = 2165.3333
$ws.Range("L132").Value = 6495.999899999999
$ws.Range("N132").Value = -11555.9999
$ws.Range("H137").Value = 1886.7858
$ws.Range("I137").Value = 1886.7858
$ws.Range("K137").Value = 5660.357400000001
$ws.Range("M137").Value = -3110.357400000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1608.3077
$ws.Range("I2").Value = 850.94446
$ws.Range("K2").Value = 850.94446
$ws.Range("M2").Value = -737.94446
$ws.Range("H32").Value = 2548.3157
$ws.Range("I32").Value = 2401.5095
$ws.Range("K32").Value = 2401.5095
$ws.Range("M32").Value = -2114.5095
$ws.Range("H45").Value = 2111
$ws.Range("I45").Value = 2249.75
$ws.Range("K45").Value = 2249.75
$ws.Range("M45").Value = -1872.75
$ws.Range("H74").Value = 3277.2222
$ws.Range("I74").Value = 3277.2222
$ws.Range("K74").Value = 3277.2222
$ws.Range("M74").Value = -2403.2222
$ws.Range("H77").Value = 3277.2222
$ws.Range("I77").Value = 3277.2222
$ws.Range("K77").Value = 16386.111
$ws.Range("M77").Value = -12018.111
$ws.Range("H110").Value = 1067.8667
$ws.Range("I110").Value = 1019.9259
$ws.Range("K110").Value = 1019.9259
$ws.Range("M110").Value = 1025.0741
$ws.Range("H116").Value = 1608.3077
$ws.Range("I116").Value = 850.94446
$ws.Range("K116").Value = 850.94446
$ws.Range("M116").Value = 1443.05554
$ws.Range("H131").Value = 103857.5
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 103857.5
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 103857.5
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -113937.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1608.3077
$ws.Range("I3").Value = 850.94446
$ws.Range("K3").Value = 850.94446
$ws.Range("M3").Value = -736.94446
$ws.Range("H20").Value = 627
$ws.Range("I20").Value = 627
$ws.Range("J20").Value = 627
$ws.Range("K20").Value = 627
$ws.Range("L20").Value = 627
$ws.Range("M20").Value = -380
$ws.Range("N20").Value = -1121
$ws.Range("H94").Value = 3709.65
$ws.Range("I94").Value = 3624.9285
$ws.Range("K94").Value = 3624.9285
$ws.Range("M94").Value = -3173.9285

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3286.4119
$ws.Range("J16").Value = 1833
$ws.Range("L16").Value = 1833
$ws.Range("N16").Value = -2407
$ws.Range("H31").Value = 7028.5454
$ws.Range("I31").Value = 3279.4614
$ws.Range("J31").Value = 12443.889
$ws.Range("K31").Value = 3279.4614
$ws.Range("L31").Value = 12443.889
$ws.Range("M31").Value = -2984.4614
$ws.Range("N31").Value = -13033.889
$ws.Range("H34").Value = 7028.5454
$ws.Range("I34").Value = 3279.4614
$ws.Range("J34").Value = 12443.889
$ws.Range("K34").Value = 3279.4614
$ws.Range("L34").Value = 12443.889
$ws.Range("M34").Value = -3077.4614
$ws.Range("N34").Value = -12847.889
$ws.Range("H107").Value = 1178.2572
$ws.Range("I107").Value = 655.6667
$ws.Range("K107").Value = 655.6667
$ws.Range("M107").Value = 1264.3333
$ws.Range("H113").Value = 3286.4119
$ws.Range("J113").Value = 1833
$ws.Range("L113").Value = 1833
$ws.Range("N113").Value = -6173

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 727.8461
$ws.Range("I7").Value = 1044.7778
$ws.Range("J7").Value = 14.75
$ws.Range("K7").Value = 3134.3334
$ws.Range("L7").Value = 44.25
$ws.Range("M7").Value = -3022.3334
$ws.Range("N7").Value = -268.25
$ws.Range("H22").Value = 357.42856
$ws.Range("J22").Value = 357.42856
$ws.Range("L22").Value = 1072.28568
$ws.Range("N22").Value = -1410.28568
$ws.Range("H23").Value = 3555.6316
$ws.Range("I23").Value = 3051.6365
$ws.Range("J23").Value = 4248.625
$ws.Range("K23").Value = 9154.9095
$ws.Range("L23").Value = 12745.875
$ws.Range("M23").Value = -8919.9095
$ws.Range("N23").Value = -13215.875
$ws.Range("H27").Value = 357.42856
$ws.Range("J27").Value = 357.42856
$ws.Range("L27").Value = 1072.28568
$ws.Range("N27").Value = -1276.28568
$ws.Range("H40").Value = 121
$ws.Range("I40").Value = 21.11111
$ws.Range("J40").Value = 270.83334
$ws.Range("K40").Value = 84.44444
$ws.Range("L40").Value = 1083.33336
$ws.Range("M40").Value = -15.44444
$ws.Range("N40").Value = -1221.33336
$ws.Range("H97").Value = 348.8421
$ws.Range("J97").Value = 386.2
$ws.Range("L97").Value = 1158.6
$ws.Range("N97").Value = -2150.6
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 867
$ws.Range("K113").Value = 2601
$ws.Range("M113").Value = -431
$ws.Range("H122").Value = 1853.2
$ws.Range("J122").Value = 2733
$ws.Range("L122").Value = 24597
$ws.Range("N122").Value = -29497

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6462.909
$ws.Range("I80").Value = 5099
$ws.Range("K80").Value = 5099
$ws.Range("M80").Value = -4101
$ws.Range("H83").Value = 6462.909
$ws.Range("I83").Value = 5099
$ws.Range("K83").Value = 25495
$ws.Range("M83").Value = -20503
$ws.Range("H102").Value = 4685.1304
$ws.Range("I102").Value = 4685.1304
$ws.Range("K102").Value = 4685.1304
$ws.Range("M102").Value = -3063.1304
$ws.Range("H113").Value = 7899.8
$ws.Range("I113").Value = 6499.8335
$ws.Range("J113").Value = 9999.75
$ws.Range("K113").Value = 6499.8335
$ws.Range("L113").Value = 9999.75
$ws.Range("M113").Value = -4329.8335
$ws.Range("N113").Value = -14339.75
$ws.Range("H122").Value = 12210.454
$ws.Range("I122").Value = 13919.471
$ws.Range("K122").Value = 41758.413
$ws.Range("M122").Value = -39308.413
$ws.Range("H132").Value = 2697.0527
$ws.Range("I132").Value = 2680.2222
$ws.Range("K132").Value = 8040.6666
$ws.Range("M132").Value = -5510.6666
$ws.Range("H141").Value = 65000
$ws.Range("J141").Value = 65000
$ws.Range("L141").Value = 65000
$ws.Range("N141").Value = -75360

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3181.2727
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 3856.2856
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 3856.2856
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -4232.2856
$ws.Range("H61").Value = 2589
$ws.Range("J61").Value = 6223.75
$ws.Range("L61").Value = 6223.75
$ws.Range("N61").Value = -6627.75
$ws.Range("H82").Value = 3875.9048
$ws.Range("I82").Value = 3207.2144
$ws.Range("J82").Value = 5213.2856
$ws.Range("K82").Value = 3207.2144
$ws.Range("L82").Value = 5213.2856
$ws.Range("M82").Value = -2846.2144
$ws.Range("N82").Value = -5935.2856
$ws.Range("H85").Value = 3875.9048
$ws.Range("I85").Value = 3207.2144
$ws.Range("J85").Value = 5213.2856
$ws.Range("K85").Value = 3207.2144
$ws.Range("L85").Value = 5213.2856
$ws.Range("M85").Value = -1959.2144
$ws.Range("N85").Value = -7709.2856
$ws.Range("H113").Value = 2589
$ws.Range("J113").Value = 6223.75
$ws.Range("L113").Value = 6223.75
$ws.Range("N113").Value = -10563.75
$ws.Range("H132").Value = 2834.0417
$ws.Range("I132").Value = 2765.7058
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 8297.117400000001
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -5767.117400000001
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 4214.033
$ws.Range("J136").Value = 3800.8823
$ws.Range("L136").Value = 11402.6469
$ws.Range("N136").Value = -16502.6469

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1599.7894
$ws.Range("I113").Value = 1020.1
$ws.Range("K113").Value = 3060.3
$ws.Range("M113").Value = -890.3000000000002
$ws.Range("H132").Value = 1116.6097
$ws.Range("I132").Value = 994
$ws.Range("J132").Value = 2669.6667
$ws.Range("K132").Value = 2982
$ws.Range("L132").Value = 8009.000100000001
$ws.Range("M132").Value = -452
$ws.Range("N132").Value = -13069.0001
$ws.Range("H141").Value = 75000
$ws.Range("J141").Value = 75000
$ws.Range("L141").Value = 75000
$ws.Range("N141").Value = -85360
